$d = $word.ActiveDocument

# Find the paragraph holding the "LOQ4205: Sistemas Produtivos II (Requisito
# fraco)" requirement line (near the end of the "Requisitos" section).
$loq = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOQ4205*") {
        $loq = $p
    }
}

if ($loq -ne $null) {
    # The footer block that follows consists of exactly three paragraphs:
    #   1) a blank "Normal" paragraph
    #   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
    #   3) the "(c) 2020 ... Creative Commons Attribution" copyright line
    # Remove all three, leaving the blank paragraph/page-break paragraph
    # that come after them untouched.
    $startPara = $loq.Next()
    $endPara = $startPara.Next().Next()

    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
